$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 21.50919633333334
$ws.Range("H2").Value = 64.52758900000001
$ws.Range("I2").Value = 0.1832723264758264
$ws.Range("J2").Value = 0.1832723264758264
$ws.Range("M2").Value = 2.148311
$ws.Range("N2").Value = 6.444933
$ws.Range("O2").Value = 0.1910770960367323
$ws.Range("P2").Value = 0.1910770960367323
$ws.Range("Q2").Value = 46.20844308405967
$ws.Range("R2").Value = 415.875987756537
$ws.Range("S2").Value = 0.03501914392689685
$ws.Range("T2").Value = 0.03501914392689683

# Row 3
$ws.Range("G3").Value = 21.50919633333334
$ws.Range("H3").Value = 64.52758900000001
$ws.Range("I3").Value = 0.1832723264758264
$ws.Range("J3").Value = 0.1832723264758264
$ws.Range("M3").Value = 4.729556333333333
$ws.Range("N3").Value = 14.188669
$ws.Range("O3").Value = 0.4206606444390354
$ws.Range("P3").Value = 0.4206606444390354
$ws.Range("Q3").Value = 101.7289557432268
$ws.Range("R3").Value = 915.5606016890411
$ws.Range("S3").Value = 0.07709545496316245
$ws.Range("T3").Value = 0.07709545496316243

# Row 4
$ws.Range("G4").Value = 21.50919633333334
$ws.Range("H4").Value = 64.52758900000001
$ws.Range("I4").Value = 0.1832723264758264
$ws.Range("J4").Value = 0.1832723264758264
$ws.Range("M4").Value = 4.365296
$ws.Range("N4").Value = 13.095888
$ws.Range("O4").Value = 0.3882622595242324
$ws.Range("P4").Value = 0.3882622595242324
$ws.Range("Q4").Value = 93.89400871711467
$ws.Range("R4").Value = 845.0460784540321
$ws.Range("S4").Value = 0.07115772758576717
$ws.Range("T4").Value = 0.07115772758576716

# Row 5
$ws.Range("I5").Value = 0.5927317426910698
$ws.Range("J5").Value = 0.5927317426910698
$ws.Range("M5").Value = 2.148311
$ws.Range("N5").Value = 6.444933
$ws.Range("O5").Value = 0.1910770960367323
$ws.Range("P5").Value = 0.1910770960367323
$ws.Range("Q5").Value = 149.445426502339
$ws.Range("R5").Value = 1345.008838521051
$ws.Range("S5").Value = 0.1132574601222012
$ws.Range("T5").Value = 0.1132574601222012

# Row 6
$ws.Range("I6").Value = 0.5927317426910698
$ws.Range("J6").Value = 0.5927317426910698
$ws.Range("M6").Value = 4.729556333333333
$ws.Range("N6").Value = 14.188669
$ws.Range("O6").Value = 0.4206606444390354
$ws.Range("P6").Value = 0.4206606444390354
$ws.Range("Q6").Value = 329.0075614758936
$ws.Range("R6").Value = 2961.068053283043
$ws.Range("S6").Value = 0.249338916859898
$ws.Range("T6").Value = 0.2493389168598979

# Row 7
$ws.Range("I7").Value = 0.5927317426910698
$ws.Range("J7").Value = 0.5927317426910698
$ws.Range("M7").Value = 4.365296
$ws.Range("N7").Value = 13.095888
$ws.Range("O7").Value = 0.3882622595242324
$ws.Range("P7").Value = 0.3882622595242324
$ws.Range("Q7").Value = 303.668101373104
$ws.Range("R7").Value = 2733.012912357936
$ws.Range("S7").Value = 0.2301353657089707
$ws.Range("T7").Value = 0.2301353657089707

# Row 8
$ws.Range("G8").Value = 26.28859766666667
$ws.Range("H8").Value = 78.86579300000001
$ws.Range("I8").Value = 0.2239959308331038
$ws.Range("J8").Value = 0.2239959308331038
$ws.Range("M8").Value = 2.148311
$ws.Range("N8").Value = 6.444933
$ws.Range("O8").Value = 0.1910770960367323
$ws.Range("P8").Value = 0.1910770960367323
$ws.Range("Q8").Value = 56.47608354187435
$ws.Range("R8").Value = 508.2847518768691
$ws.Range("S8").Value = 0.04280049198763421
$ws.Range("T8").Value = 0.04280049198763421

# Row 9
$ws.Range("G9").Value = 26.28859766666667
$ws.Range("H9").Value = 78.86579300000001
$ws.Range("I9").Value = 0.2239959308331038
$ws.Range("J9").Value = 0.2239959308331038
$ws.Range("M9").Value = 4.729556333333333
$ws.Range("N9").Value = 14.188669
$ws.Range("O9").Value = 0.4206606444390354
$ws.Range("P9").Value = 0.4206606444390354
$ws.Range("Q9").Value = 124.3334035888352
$ws.Range("R9").Value = 1119.000632299517
$ws.Range("S9").Value = 0.09422627261597505
$ws.Range("T9").Value = 0.09422627261597502

# Row 10
$ws.Range("G10").Value = 26.28859766666667
$ws.Range("H10").Value = 78.86579300000001
$ws.Range("I10").Value = 0.2239959308331038
$ws.Range("J10").Value = 0.2239959308331038
$ws.Range("M10").Value = 4.365296
$ws.Range("N10").Value = 13.095888
$ws.Range("O10").Value = 0.3882622595242324
$ws.Range("P10").Value = 0.3882622595242324
$ws.Range("Q10").Value = 114.7575102399094
$ws.Range("R10").Value = 1032.817592159184
$ws.Range("S10").Value = 0.08696916622949455
$ws.Range("T10").Value = 0.08696916622949454

